$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the blank password placeholder (a single space) from cell A2.
# This removes the now-unused shared string " " from the shared string
# table, shifting all later shared-string indices down by one, while the
# remaining rows (A3:A19) keep their original row numbers and values.
$ws.Range("A2").ClearContents()

# Narrow the active selection down to just A2 (matches the saved view state).
$ws.Range("A2").Select()
